$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report values between row 2 and row 4
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
#  Origen, Precio $/Kg) while leaving the rest of each row untouched.

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")
    $temp = $cell2.Value2
    $cell2.Value2 = $cell4.Value2
    $cell4.Value2 = $temp
}
